$wb = $excel.ActiveWorkbook

# Sheet "Hoja1" holds the daily conversion note in A1.
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.36 = 12733.02 pesos`n✅ 12733.02 pesos = 3.33 = 966.49 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas" holds the rate table with updated N10/O10/N12/O12 values.
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 298
$ws2.Range("O10").Value = 3794.44
$ws2.Range("N12").Value = 3826
$ws2.Range("O12").Value = 290.41
